$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Restore C10 to its prior value (18 -> 1) per the target revision.
$ws.Range("C10").Value = 1
